$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the ID block (old row 55), pushing
# everything from the old row 55 onward down by two rows.
$ws.Rows("55:56").Insert()

# Seed the two new rows with formatting copied from an existing
# "age group" row (row 2: style s="5" for col A, s="3" for col B/C).
$ws.Range("A2:C2").Copy($ws.Range("A55"))
$ws.Range("A2:C2").Copy($ws.Range("A56"))

# Column C on the age-group template rows normally carries no explicit
# style, but the two new rows need C formatted with style index 3
# (same as column B). Pull that exact style from an existing cell that
# already uses it.
$ws.Range("C43").Copy()
$ws.Range("C55").PasteSpecial(-4122)
$ws.Range("C56").PasteSpecial(-4122)

# Fill in the new participant IDs and their age values.
$ws.Range("A55").Value2 = "blg077"
$ws.Range("B55").Value2 = 5
$ws.Range("C55").Value2 = 5

$ws.Range("A56").Value2 = "blg085"
$ws.Range("B56").Value2 = 5
$ws.Range("C56").Value2 = 5

# Re-apply the sort so the sortState/sortCondition refs cover the new
# full range of data (A2:C97 / A2:A97).
$s = $ws.Sort
$s.SortFields.Add($ws.Range("A2:A97"))
$s.SetRange($ws.Range("A2:C97"))
$s.Header = -4142
$s.Apply()

# Match the author's final view state: scrolled back to the top with
# C57 selected.
$ws.Range("C57").Select()
